$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.08805033333333334
$ws.Range("H2").Value = 0.264151
$ws.Range("I2").Value = 0.005589762818257384
$ws.Range("J2").Value = 0.005589762818257385
$ws.Range("M2").Value = 127.3992563333333
$ws.Range("N2").Value = 382.197769
$ws.Range("O2").Value = 0.4838549810199306
$ws.Range("P2").Value = 0.4838549810199307
$ws.Range("Q2").Value = 11.21754698656878
$ws.Range("R2").Value = 100.957922879119
$ws.Range("S2").Value = 0.002704634582333841
$ws.Range("T2").Value = 0.002704634582333842

$ws.Range("G3").Value = 0.08805033333333334
$ws.Range("H3").Value = 0.264151
$ws.Range("I3").Value = 0.005589762818257384
$ws.Range("J3").Value = 0.005589762818257385
$ws.Range("M3").Value = 59.36586533333332
$ws.Range("O3").Value = 0.2254681108101269
$ws.Range("P3").Value = 0.2254681108101269
$ws.Range("Q3").Value = 5.227184231221777
$ws.Range("R3").Value = 47.044658080996
$ws.Range("S3").Value = 0.001260313262509183
$ws.Range("T3").Value = 0.001260313262509183

$ws.Range("G4").Value = 0.08805033333333334
$ws.Range("H4").Value = 0.264151
$ws.Range("I4").Value = 0.005589762818257384
$ws.Range("J4").Value = 0.005589762818257385
$ws.Range("M4").Value = 16.63275166666667
$ws.Range("N4").Value = 49.898255
$ws.Range("O4").Value = 0.06317022542837675
$ws.Range("P4").Value = 0.06317022542837675
$ws.Range("Q4").Value = 1.464519328500556
$ws.Range("R4").Value = 13.180673956505
$ws.Range("S4").Value = 0.0003531065773204775
$ws.Range("T4").Value = 0.0003531065773204776

$ws.Range("G5").Value = 0.08805033333333334
$ws.Range("H5").Value = 0.264151
$ws.Range("I5").Value = 0.005589762818257384
$ws.Range("J5").Value = 0.005589762818257385
$ws.Range("M5").Value = 59.90262233333334
$ws.Range("N5").Value = 179.707867
$ws.Range("O5").Value = 0.2275066827415657
$ws.Range("P5").Value = 0.2275066827415658
$ws.Range("Q5").Value = 5.274445863990779
$ws.Range("R5").Value = 47.47001277591701
$ws.Range("S5").Value = 0.001271708396093883
$ws.Range("T5").Value = 0.001271708396093884

$ws.Range("I6").Value = 0.9470512964761942
$ws.Range("J6").Value = 0.9470512964761943
$ws.Range("M6").Value = 127.3992563333333
$ws.Range("N6").Value = 382.197769
$ws.Range("O6").Value = 0.4838549810199306
$ws.Range("P6").Value = 0.4838549810199307
$ws.Range("Q6").Value = 1900.544399167281
$ws.Range("R6").Value = 17104.89959250553
$ws.Range("S6").Value = 0.4582354870813897
$ws.Range("T6").Value = 0.4582354870813898

$ws.Range("I7").Value = 0.9470512964761942
$ws.Range("J7").Value = 0.9470512964761943
$ws.Range("M7").Value = 59.36586533333332
$ws.Range("O7").Value = 0.2254681108101269
$ws.Range("P7").Value = 0.2254681108101269
$ws.Range("S7").Value = 0.2135298666567689
$ws.Range("T7").Value = 0.2135298666567689

$ws.Range("I8").Value = 0.9470512964761942
$ws.Range("J8").Value = 0.9470512964761943
$ws.Range("M8").Value = 16.63275166666667
$ws.Range("N8").Value = 49.898255
$ws.Range("O8").Value = 0.06317022542837675
$ws.Range("P8").Value = 0.06317022542837675
$ws.Range("Q8").Value = 248.1276887528634
$ws.Range("R8").Value = 2233.14919877577
$ws.Range("S8").Value = 0.05982544389063765
$ws.Range("T8").Value = 0.05982544389063766

$ws.Range("I9").Value = 0.9470512964761942
$ws.Range("J9").Value = 0.9470512964761943
$ws.Range("M9").Value = 59.90262233333334
$ws.Range("N9").Value = 179.707867
$ws.Range("O9").Value = 0.2275066827415657
$ws.Range("P9").Value = 0.2275066827415658
$ws.Range("Q9").Value = 893.6283982158689
$ws.Range("R9").Value = 8042.655583942819
$ws.Range("S9").Value = 0.215460498847398
$ws.Range("T9").Value = 0.2154604988473981

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.7460013333333334
$ws.Range("H10").Value = 2.238004
$ws.Range("I10").Value = 0.04735894070554834
$ws.Range("J10").Value = 0.04735894070554835
$ws.Range("M10").Value = 127.3992563333333
$ws.Range("N10").Value = 382.197769
$ws.Range("O10").Value = 0.4838549810199306
$ws.Range("P10").Value = 0.4838549810199307
$ws.Range("Q10").Value = 95.04001509034178
$ws.Range("R10").Value = 855.360135813076
$ws.Range("S10").Value = 0.02291485935620712
$ws.Range("T10").Value = 0.02291485935620712

$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.7460013333333334
$ws.Range("H11").Value = 2.238004
$ws.Range("I11").Value = 0.04735894070554834
$ws.Range("J11").Value = 0.04735894070554835
$ws.Range("M11").Value = 59.36586533333332
$ws.Range("O11").Value = 0.2254681108101269
$ws.Range("P11").Value = 0.2254681108101269
$ws.Range("Q11").Value = 44.28701469315378
$ws.Range("R11").Value = 398.583132238384
$ws.Range("S11").Value = 0.0106779308908488
$ws.Range("T11").Value = 0.0106779308908488

$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.7460013333333334
$ws.Range("H12").Value = 2.238004
$ws.Range("I12").Value = 0.04735894070554834
$ws.Range("J12").Value = 0.04735894070554835
$ws.Range("M12").Value = 16.63275166666667
$ws.Range("N12").Value = 49.898255
$ws.Range("O12").Value = 0.06317022542837675
$ws.Range("P12").Value = 0.06317022542837675
$ws.Range("Q12").Value = 12.40805492033556
$ws.Range("R12").Value = 111.67249428302
$ws.Range("S12").Value = 0.002991674960418617
$ws.Range("T12").Value = 0.002991674960418617

$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.7460013333333334
$ws.Range("H13").Value = 2.238004
$ws.Range("I13").Value = 0.04735894070554834
$ws.Range("J13").Value = 0.04735894070554835
$ws.Range("M13").Value = 59.90262233333334
$ws.Range("N13").Value = 179.707867
$ws.Range("O13").Value = 0.2275066827415657
$ws.Range("P13").Value = 0.2275066827415658
$ws.Range("Q13").Value = 44.68743613082979
$ws.Range("R13").Value = 402.1869251774681
$ws.Range("S13").Value = 0.01077447549807381
$ws.Range("T13").Value = 0.01077447549807382
